$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 and 4 entirely (trailing data removed)
$ws.Rows("3:4").Delete()

# Update the remaining label values
$ws.Range("A1").Value = "A"
$ws.Range("A2").Value = "B"

# Update selection to match final cursor position
$ws.Range("B2").Select()
